{"js": "// Word JS API (Office.js) script.\n// Body of: async (context) => { ... }\n//\n// Change 1: the Title paragraph's text (\"Modern Dive Chapter 2: Data\n//           Visualization\") was split across many single-word runs;\n//           consolidate it into one run (same text, same paragraph style).\n// Change 2: the built-in \"Subtitle\" paragraph style now inherits from\n//           \"Title\" instead of \"Normal\", and no longer carries an explicit\n//           (theme) font color override.\n// Change 3: the \"Abstract Title\" paragraph style no longer carries an\n//           explicit font color override.\n\n// --- Change 1: collapse the Title paragraph's runs into a single run ---\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"style,text\"));\nawait context.sync();\n\nconst titlePara = paragraphs.items.find((p) => p.style === \"Title\");\nif (titlePara) {\n  const fullText = titlePara.text;\n  titlePara.getRange().insertText(fullText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2 & 3: update paragraph style definitions ---\nconst styles = context.document.getStyles();\nstyles.load(\"items\");\nawait context.sync();\n\nconst subtitleStyle = styles.getByNameOrNullObject(\"Subtitle\");\nconst abstractTitleStyle = styles.getByNameOrNullObject(\"Abstract Title\");\nawait context.sync();\n\nsubtitleStyle.load(\"isNullObject\");\nabstractTitleStyle.load(\"isNullObject\");\nawait context.sync();\n\nif (!subtitleStyle.isNullObject) {\n  subtitleStyle.baseStyle = \"Title\";\n  // Clear the explicit (theme) font color override -> falls back to\n  // automatic, same as the now-removed <w:color> element.\n  subtitleStyle.font.color = -16777216 /* wdColorAutomatic */;\n}\n\nif (!abstractTitleStyle.isNullObject) {\n  abstractTitleStyle.font.color = -16777216 /* wdColorAutomatic */;\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop script (PowerShell-style).\n# $word.ActiveDocument is the open document.\n#\n# Change 1: the Title paragraph's text (\"Modern Dive Chapter 2: Data\n#           Visualization\") was split across many single-word runs;\n#           consolidate it into one run (same text, same paragraph style).\n# Change 2: the built-in \"Subtitle\" paragraph style now inherits from\n#           \"Title\" instead of \"Normal\", and no longer carries an explicit\n#           (theme) font color override.\n# Change 3: the \"Abstract Title\" paragraph style no longer carries an\n#           explicit font color override.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: collapse the Title paragraph's runs into a single run ---\n$titlePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Style.NameLocal -eq \"Title\") {\n        $titlePara = $p\n        break\n    }\n}\n\nif ($titlePara -ne $null) {\n    $titleText = $titlePara.Range.Text\n    # Paragraph.Range.Text includes the trailing paragraph mark; strip it\n    # before round-tripping the text through Find/Replace.\n    $titleText = $titleText.TrimEnd([char]13, [char]7)\n\n    $find = $d.Content.Find\n    $find.Text = $titleText\n    $find.Replacement.Text = $titleText\n    $find.Execute($titleText, $false, $false, $false, $false, $false, $true, 1, $false, $titleText, 2)\n}\n\n# --- Change 2 & 3: update paragraph style definitions ---\n$subtitleStyle = $d.Styles(\"Subtitle\")\n$subtitleStyle.BaseStyle = $d.Styles(\"Title\")\n$subtitleStyle.Font.Color = -16777216   # wdColorAutomatic\n\n$abstractTitleStyle = $d.Styles(\"Abstract Title\")\n$abstractTitleStyle.Font.Color = -16777216   # wdColorAutomatic\n"}
